# Append a new data row (row 94) to the sheet, mirroring the existing
# rows in data/KK.MI.xlsx (results produced by the R script).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 93
$newRow = 94

# Copy the formatting (styles) of the previous row down into the new
# row so the date column keeps its date/time number format, etc.
$ws.Range("A" + $lastRow + ":H" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":H" + $newRow).PasteSpecial(-4122)  # xlPasteFormats

# date (2024-06-06 07:00:00, serial 45449.2916666667)
$ws.Cells.Item($newRow, 1).Value = 45449.2916666667
# volume
$ws.Cells.Item($newRow, 2).Value = 0
# high
$ws.Cells.Item($newRow, 3).Value = 2
# low
$ws.Cells.Item($newRow, 4).Value = 2
# open
$ws.Cells.Item($newRow, 5).Value = 2
# close
$ws.Cells.Item($newRow, 6).Value = 2

# adj_close is stored as text "2" (shared string), not a number -
# force text formatting before assigning so it isn't reinterpreted
# as a number, then drop back to the default style (no explicit
# number format on the cell itself).
$ws.Range("G" + $newRow).NumberFormat = "@"
$ws.Range("G" + $newRow).Value = "2"
$ws.Range("G" + $newRow).Style = "Normal"

# ticker
$ws.Cells.Item($newRow, 8).Value = "KK.MI"
